$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 202, shifting existing rows 202:220 down to 203:221
$ws.Rows(202).Insert()

# Populate the new row 202 with the new record
$ws.Range("A202").Value = 4
$ws.Range("B202").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C202").Value = "Los Lagos"
$ws.Range("D202").NumberFormat = $ws.Range("D201").NumberFormat
$ws.Range("D202").Value = 44578
$ws.Range("E202").Value = 10
$ws.Range("F202").Value = 100112037
$ws.Range("G202").Value = "Cebollín"
$ws.Range("H202").Value = "Sin especificar"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 60
$ws.Range("K202").Value = 6000
$ws.Range("L202").Value = 6000
$ws.Range("M202").Value = 6000
$ws.Range("N202").Value = "$/paquete 36 unidades"
$ws.Range("O202").Value = "Región Metropolitana"
$ws.Range("P202").Value = 167
$ws.Range("Q202").Value = 36
$ws.Range("R202").Value = "Hortaliza"
